$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Efna5"
$arr[0,2] = "Epha3"
$arr[0,3] = "ECs"
$arr[0,4] = 2.0
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.3227736666666667
$arr[0,7] = 0.968321
$arr[0,8] = 0.1416094457286952
$arr[0,9] = 0.1416094457286952
$arr[0,10] = 2.0
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.118034
$arr[0,13] = 3.354102
$arr[0,14] = 0.01817106018145251
$arr[0,15] = 0.01817106018145251
$arr[0,16] = 0.360871933638
$arr[0,17] = 3.247847402742
$arr[0,18] = 0.002573193760598252
$arr[0,19] = 0.002573193760598253

$arr[1,0] = "ECs"
$arr[1,1] = "Efna5"
$arr[1,2] = "Epha3"
$arr[1,3] = "FAPs"
$arr[1,4] = 2.0
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.3227736666666667
$arr[1,7] = 0.968321
$arr[1,8] = 0.1416094457286952
$arr[1,9] = 0.1416094457286952
$arr[1,10] = 3.0
$arr[1,11] = 1.0
$arr[1,12] = 59.24481466666666
$arr[1,13] = 177.734444
$arr[1,14] = 0.962887615892719
$arr[1,15] = 0.9628876158927191
$arr[1,16] = 19.12266606094711
$arr[1,17] = 172.103994548524
$arr[1,18] = 0.1363539815855927
$arr[1,19] = 0.1363539815855927

$arr[2,0] = "ECs"
$arr[2,1] = "Efna5"
$arr[2,2] = "Epha3"
$arr[2,3] = "MuSCs"
$arr[2,4] = 2.0
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.3227736666666667
$arr[2,7] = 0.968321
$arr[2,8] = 0.1416094457286952
$arr[2,9] = 0.1416094457286952
$arr[2,10] = 3.0
$arr[2,11] = 1.0
$arr[2,12] = 1.121724666666667
$arr[2,13] = 3.365174
$arr[2,14] = 0.01823104344324033
$arr[2,15] = 0.01823104344324033
$arr[2,16] = 0.3620631836504445
$arr[2,17] = 3.258568652854
$arr[2,18] = 0.002581687957053025
$arr[2,19] = 0.002581687957053025

$arr[3,0] = "ECs"
$arr[3,1] = "Efna5"
$arr[3,2] = "Epha3"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 2.0
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.3227736666666667
$arr[3,7] = 0.968321
$arr[3,8] = 0.1416094457286952
$arr[3,9] = 0.1416094457286952
$arr[3,10] = 1.0
$arr[3,11] = 0.3333333333333333
$arr[3,12] = 0.04370233333333334
$arr[3,13] = 0.131107
$arr[3,14] = 0.0007102804825880949
$arr[3,15] = 0.0007102804825880949
$arr[3,16] = 0.01410596237188889
$arr[3,17] = 0.126953661347
$arr[3,18] = 0.0001005824254512102
$arr[3,19] = 0.0001005824254512102

$arr[4,0] = "FAPs"
$arr[4,1] = "Efna5"
$arr[4,2] = "Epha3"
$arr[4,3] = "ECs"
$arr[4,4] = 3.0
$arr[4,5] = 1.0
$arr[4,6] = 1.874986333333333
$arr[4,7] = 5.624959
$arr[4,8] = 0.8226066833587575
$arr[4,9] = 0.8226066833587576
$arr[4,10] = 2.0
$arr[4,11] = 0.6666666666666666
$arr[4,12] = 1.118034
$arr[4,13] = 3.354102
$arr[4,14] = 0.01817106018145251
$arr[4,15] = 0.01817106018145251
$arr[4,16] = 2.096298470202
$arr[4,17] = 18.866686231818
$arr[4,18] = 0.01494763554897703
$arr[4,19] = 0.01494763554897704

$arr[5,0] = "FAPs"
$arr[5,1] = "Efna5"
$arr[5,2] = "Epha3"
$arr[5,3] = "FAPs"
$arr[5,4] = 3.0
$arr[5,5] = 1.0
$arr[5,6] = 1.874986333333333
$arr[5,7] = 5.624959
$arr[5,8] = 0.8226066833587575
$arr[5,9] = 0.8226066833587576
$arr[5,10] = 3.0
$arr[5,11] = 1.0
$arr[5,12] = 59.24481466666666
$arr[5,13] = 177.734444
$arr[5,14] = 0.962887615892719
$arr[5,15] = 0.9628876158927191
$arr[5,16] = 111.0832178208662
$arr[5,17] = 999.748960387796
$arr[5,18] = 0.7920777881567308
$arr[5,19] = 0.792077788156731

$arr[6,0] = "FAPs"
$arr[6,1] = "Efna5"
$arr[6,2] = "Epha3"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3.0
$arr[6,5] = 1.0
$arr[6,6] = 1.874986333333333
$arr[6,7] = 5.624959
$arr[6,8] = 0.8226066833587575
$arr[6,9] = 0.8226066833587576
$arr[6,10] = 3.0
$arr[6,11] = 1.0
$arr[6,12] = 1.121724666666667
$arr[6,13] = 3.365174
$arr[6,14] = 0.01823104344324033
$arr[6,15] = 0.01823104344324033
$arr[6,16] = 2.103218419762889
$arr[6,17] = 18.928965777866
$arr[6,18] = 0.01499697818101335
$arr[6,19] = 0.01499697818101335

$arr[7,0] = "FAPs"
$arr[7,1] = "Efna5"
$arr[7,2] = "Epha3"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3.0
$arr[7,5] = 1.0
$arr[7,6] = 1.874986333333333
$arr[7,7] = 5.624959
$arr[7,8] = 0.8226066833587575
$arr[7,9] = 0.8226066833587576
$arr[7,10] = 1.0
$arr[7,11] = 0.3333333333333333
$arr[7,12] = 0.04370233333333334
$arr[7,13] = 0.131107
$arr[7,14] = 0.0007102804825880949
$arr[7,15] = 0.0007102804825880949
$arr[7,16] = 0.08194127773477779
$arr[7,17] = 0.7374714996130001
$arr[7,18] = 0.0005842814720362504
$arr[7,19] = 0.0005842814720362506

$arr[8,0] = "MuSCs"
$arr[8,1] = "Efna5"
$arr[8,2] = "Epha3"
$arr[8,3] = "ECs"
$arr[8,4] = 2.0
$arr[8,5] = 0.6666666666666666
$arr[8,6] = 0.081563
$arr[8,7] = 0.244689
$arr[8,8] = 0.03578387091254728
$arr[8,9] = 0.03578387091254728
$arr[8,10] = 2.0
$arr[8,11] = 0.6666666666666666
$arr[8,12] = 1.118034
$arr[8,13] = 3.354102
$arr[8,14] = 0.01817106018145251
$arr[8,15] = 0.01817106018145251
$arr[8,16] = 0.09119020714199999
$arr[8,17] = 0.820711864278
$arr[8,18] = 0.0006502308718772244
$arr[8,19] = 0.0006502308718772246

$arr[9,0] = "MuSCs"
$arr[9,1] = "Efna5"
$arr[9,2] = "Epha3"
$arr[9,3] = "FAPs"
$arr[9,4] = 2.0
$arr[9,5] = 0.6666666666666666
$arr[9,6] = 0.081563
$arr[9,7] = 0.244689
$arr[9,8] = 0.03578387091254728
$arr[9,9] = 0.03578387091254728
$arr[9,10] = 3.0
$arr[9,11] = 1.0
$arr[9,12] = 59.24481466666666
$arr[9,13] = 177.734444
$arr[9,14] = 0.962887615892719
$arr[9,15] = 0.9628876158927191
$arr[9,16] = 4.832184818657333
$arr[9,17] = 43.48966336791599
$arr[9,18] = 0.03445584615039546
$arr[9,19] = 0.03445584615039547

$arr[10,0] = "MuSCs"
$arr[10,1] = "Efna5"
$arr[10,2] = "Epha3"
$arr[10,3] = "MuSCs"
$arr[10,4] = 2.0
$arr[10,5] = 0.6666666666666666
$arr[10,6] = 0.081563
$arr[10,7] = 0.244689
$arr[10,8] = 0.03578387091254728
$arr[10,9] = 0.03578387091254728
$arr[10,10] = 3.0
$arr[10,11] = 1.0
$arr[10,12] = 1.121724666666667
$arr[10,13] = 3.365174
$arr[10,14] = 0.01823104344324033
$arr[10,15] = 0.01823104344324033
$arr[10,16] = 0.09149122898733333
$arr[10,17] = 0.823421060886
$arr[10,18] = 0.0006523773051739534
$arr[10,19] = 0.0006523773051739534

$arr[11,0] = "MuSCs"
$arr[11,1] = "Efna5"
$arr[11,2] = "Epha3"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 2.0
$arr[11,5] = 0.6666666666666666
$arr[11,6] = 0.081563
$arr[11,7] = 0.244689
$arr[11,8] = 0.03578387091254728
$arr[11,9] = 0.03578387091254728
$arr[11,10] = 1.0
$arr[11,11] = 0.3333333333333333
$arr[11,12] = 0.04370233333333334
$arr[11,13] = 0.131107
$arr[11,14] = 0.0007102804825880949
$arr[11,15] = 0.0007102804825880949
$arr[11,16] = 0.003564493413666667
$arr[11,17] = 0.032080440723
$arr[11,18] = 0.00002541658510063417
$arr[11,19] = 0.00002541658510063417

$ws.Range("A2:T13").Value = $arr